$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 15: hours, release and description for a new log entry.
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = "1.5a"
$ws.Range("D15").Value = "Attachments & stability."

# Move the active selection to D16, as in the authored workbook.
$ws.Range("D16").Select()

# Recalculate so the SUM formula in C24 reflects the new hours.
$excel.Calculate()
